$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2170
$ws.Range("I33").Value = 1845.6666
$ws.Range("J33").Value = 3467.3333
$ws.Range("K33").Value = 1845.6666
$ws.Range("L33").Value = 3467.3333
$ws.Range("M33").Value = -1616.6666
$ws.Range("N33").Value = -3925.3333
# Row 86
$ws.Range("H86").Value = 14100.223
$ws.Range("I86").Value = 2501.5
$ws.Range("J86").Value = 17414.143
$ws.Range("K86").Value = 2501.5
$ws.Range("L86").Value = 17414.143
$ws.Range("M86").Value = -1378.5
$ws.Range("N86").Value = -19660.143
# Row 89
$ws.Range("H89").Value = 14100.223
$ws.Range("I89").Value = 2501.5
$ws.Range("J89").Value = 17414.143
$ws.Range("K89").Value = 12507.5
$ws.Range("L89").Value = 87070.715
$ws.Range("M89").Value = -6891.5
$ws.Range("N89").Value = -98302.715
# Row 100
$ws.Range("H100").Value = 25002426
$ws.Range("I100").Value = 1133.6666
$ws.Range("J100").Value = 40003200
$ws.Range("K100").Value = 1133.6666
$ws.Range("L100").Value = 40003200
$ws.Range("M100").Value = -592.6666
$ws.Range("N100").Value = -40004282
# Row 113
$ws.Range("H113").Value = 4001756.5
$ws.Range("I113").Value = 5883971
$ws.Range("J113").Value = 2050
$ws.Range("K113").Value = 5883971
$ws.Range("L113").Value = 2050
$ws.Range("M113").Value = -5880717
$ws.Range("N113").Value = -8558
# Row 116
$ws.Range("H116").Value = 16673988
$ws.Range("I116").Value = 8335558.5
$ws.Range("J116").Value = 27791896
$ws.Range("K116").Value = 8335558.5
$ws.Range("L116").Value = 27791896
$ws.Range("M116").Value = -8332116.5
$ws.Range("N116").Value = -27798780

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 18654
$ws.Range("I2").Value = 35883
$ws.Range("J2").Value = 1425
$ws.Range("K2").Value = 35883
$ws.Range("L2").Value = 1425
$ws.Range("M2").Value = -35770
$ws.Range("N2").Value = -1651
# Row 5
$ws.Range("H5").Value = 993.3333
$ws.Range("I5").Value = 480
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 480
$ws.Range("L5").Value = 1250
$ws.Range("M5").Value = -368
$ws.Range("N5").Value = -1474
# Row 45
$ws.Range("H45").Value = 417814.28
$ws.Range("I45").Value = 1000995
$ws.Range("J45").Value = 1256.6428
$ws.Range("K45").Value = 1000995
$ws.Range("L45").Value = 1256.6428
$ws.Range("M45").Value = -1000618
$ws.Range("N45").Value = -2010.6428
# Row 61
$ws.Range("H61").Value = 3466301.8
$ws.Range("I61").Value = 1895164.4
$ws.Range("J61").Value = 8404162
$ws.Range("K61").Value = 1895164.4
$ws.Range("L61").Value = 8404162
$ws.Range("M61").Value = -1894952.4
$ws.Range("N61").Value = -8404586
# Row 102
$ws.Range("H102").Value = 2581.818
$ws.Range("I102").Value = 2775
$ws.Range("J102").Value = 2066.6667
$ws.Range("K102").Value = 2775
$ws.Range("L102").Value = 2066.6667
$ws.Range("M102").Value = -1153
$ws.Range("N102").Value = -5310.6667
# Row 110
$ws.Range("H110").Value = 719.1818
$ws.Range("I110").Value = 816.7143
$ws.Range("J110").Value = 548.5
$ws.Range("K110").Value = 816.7143
$ws.Range("L110").Value = 548.5
$ws.Range("M110").Value = 1228.2857
$ws.Range("N110").Value = -4638.5
# Row 116
$ws.Range("H116").Value = 18654
$ws.Range("I116").Value = 35883
$ws.Range("J116").Value = 1425
$ws.Range("K116").Value = 35883
$ws.Range("L116").Value = 1425
$ws.Range("M116").Value = -33589
$ws.Range("N116").Value = -6013
# Row 136
$ws.Range("H136").Value = 3466301.8
$ws.Range("I136").Value = 1895164.4
$ws.Range("J136").Value = 8404162
$ws.Range("K136").Value = 5685493.199999999
$ws.Range("L136").Value = 25212486
$ws.Range("M136").Value = -5682943.199999999
$ws.Range("N136").Value = -25217586

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 18654
$ws.Range("I3").Value = 35883
$ws.Range("J3").Value = 1425
$ws.Range("K3").Value = 35883
$ws.Range("L3").Value = 1425
$ws.Range("M3").Value = -35769
$ws.Range("N3").Value = -1653
# Row 4
$ws.Range("H4").Value = 993.3333
$ws.Range("I4").Value = 480
$ws.Range("J4").Value = 1250
$ws.Range("K4").Value = 480
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = -365
$ws.Range("N4").Value = -1480
# Row 20
$ws.Range("H20").Value = 20012262
$ws.Range("I20").Value = 41677132
$ws.Range("J20").Value = 13920.615
$ws.Range("K20").Value = 41677132
$ws.Range("L20").Value = 13920.615
$ws.Range("M20").Value = -41676885
$ws.Range("N20").Value = -14414.615
# Row 22
$ws.Range("H22").Value = 482.91666
$ws.Range("I22").Value = 473.3913
$ws.Range("J22").Value = 702
$ws.Range("K22").Value = 473.3913
$ws.Range("L22").Value = 702
$ws.Range("M22").Value = -300.3913
$ws.Range("N22").Value = -1048
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
# Row 94
$ws.Range("H94").Value = 1969.84
$ws.Range("I94").Value = 1655.2307
$ws.Range("J94").Value = 2310.6667
$ws.Range("K94").Value = 1655.2307
$ws.Range("L94").Value = 2310.6667
$ws.Range("M94").Value = -1204.2307
$ws.Range("N94").Value = -3212.6667
# Row 99
$ws.Range("H99").Value = 1999.75
$ws.Range("I99").Value = 1975
$ws.Range("J99").Value = 2012.125
$ws.Range("K99").Value = 1975
$ws.Range("L99").Value = 2012.125
$ws.Range("M99").Value = -477
$ws.Range("N99").Value = -5008.125
# Row 105
$ws.Range("H105").Value = 1766.6666
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -4794

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 315.15384
$ws.Range("I22").Value = 291.41666
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 291.41666
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 58.58334000000002
$ws.Range("N22").Value = -1300

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1257.5106
$ws.Range("I113").Value = 513.1739
$ws.Range("J113").Value = 1970.8334
$ws.Range("K113").Value = 1539.5217
$ws.Range("L113").Value = 5912.5002
$ws.Range("M113").Value = 630.4783
$ws.Range("N113").Value = -10252.5002
# Row 114
$ws.Range("H114").Value = 2161
$ws.Range("I114").Value = 620
$ws.Range("J114").Value = 2721.3635
$ws.Range("K114").Value = 1860
$ws.Range("L114").Value = 8164.0905
$ws.Range("M114").Value = 1394
$ws.Range("N114").Value = -14672.0905

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 22729836
$ws.Range("I97").Value = 1655.7142
$ws.Range("J97").Value = 62504150
$ws.Range("K97").Value = 1655.7142
$ws.Range("L97").Value = 62504150
$ws.Range("M97").Value = -1159.7142
$ws.Range("N97").Value = -62505142
# Row 111
$ws.Range("H111").Value = 30000.666
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 30000.666
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 30000.666
$ws.Range("N111").Value = -36134.666

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 35719828
$ws.Range("I22").Value = 4600
$ws.Range("J22").Value = 55561624
$ws.Range("K22").Value = 4600
$ws.Range("L22").Value = 55561624
$ws.Range("M22").Value = -4305
$ws.Range("N22").Value = -55562214
# Row 27
$ws.Range("H27").Value = 35719828
$ws.Range("I27").Value = 4600
$ws.Range("J27").Value = 55561624
$ws.Range("K27").Value = 4600
$ws.Range("L27").Value = 55561624
$ws.Range("M27").Value = -4493
$ws.Range("N27").Value = -55561838
# Row 132
$ws.Range("H132").Value = 1451016.6
$ws.Range("I132").Value = 1853104.5
$ws.Range("J132").Value = 3500.1
$ws.Range("K132").Value = 5559313.5
$ws.Range("L132").Value = 10500.3
$ws.Range("M132").Value = -5556783.5
$ws.Range("N132").Value = -15560.3
